$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $headerCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value2
    if ($text -ne $null) {
        $newText = $text -replace "_old$", "_FV2310"
        $newText = $newText -replace "_new$", "_FV2404"
        if ($newText -ne $text) {
            $cell.Value2 = $newText
        }
    }
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
